# "Big fix when updating data"
# Updates the example dates / flags / intervals used throughout the
# "main", "search" and "week_week" sheets, and refreshes the current
# selections to match the new working cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "main"
# ---------------------------------------------------------------------
$main = $wb.Worksheets.Item("main")

# Inital Date / Final Date example values (stored as TEXT, per the
# "Column C5 AND D5 have to be as TEXT data" note in I4)
$main.Range("C5").Value = "2015-05-01"
$main.Range("D5").Value = "2015-05-30"

# UPDATE ROWS flags for aFRR_Energy / mFRR_Energy / aFRR_power switched off
$main.Range("D11").Value = $false
$main.Range("E11").Value = $false
$main.Range("D12").Value = $false
$main.Range("E12").Value = $false
$main.Range("D13").Value = $false
$main.Range("E13").Value = $false

# refresh selection (was the merged E4:F5, now a single cell below it)
$main.Range("E6").Select()

# ---------------------------------------------------------------------
# Sheet "search"
# ---------------------------------------------------------------------
$search = $wb.Worksheets.Item("search")

$search.Range("C5").Value = "2015-05-01"
$search.Range("D5").Value = "2015-05-30"
$search.Range("G5").Value = "RR"

# the trailing, unused H4 cell is removed entirely (dimension C4:H5 -> C4:G5)
$search.Range("H4").Clear()

# refresh selection
$search.Range("G5").Select()

# ---------------------------------------------------------------------
# Sheet "week_week"
# ---------------------------------------------------------------------
$week = $wb.Worksheets.Item("week_week")

$week.Range("C5").Value = "2021-07-12"
$week.Range("D5").Value = 469

$week.Range("F11").Value = 60
$week.Range("F12").Value = 60
$week.Range("F13").Value = 60

# keep selection/active sheet on week_week, matching the unchanged activeTab
$week.Range("C5").Select()
